$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row right after the header (new weekly price observation),
# pushing the existing rows 6..34 down to 7..35.
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = 3
$ws.Range("B6").Value = "Femacal de La Calera"
$ws.Range("C6").Value = "Coquimbo"
$ws.Range("D6").Value = 45092
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 100112044
$ws.Range("G6").Value = "Perejil"
$ws.Range("H6").Value = "Sin especificar"
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 90
$ws.Range("K6").Value = 3000
$ws.Range("L6").Value = 3500
$ws.Range("M6").Value = 3278
$ws.Range("N6").Value = "`$/docena de atados (3 kilos)"
$ws.Range("O6").Value = "Provincia de Quillota"
$ws.Range("P6").Value = 1093
$ws.Range("Q6").Value = 3
$ws.Range("R6").Value = "Hortaliza"
